# Applies the "review-comment" colouring pass described in the commit
# message "all corrections done apart from references": every bullet in
# the corrections list that has already been actioned gets its text
# turned green (RGB 00CC33), except the bullet about the sample/pulse
# axes (P142-144), which additionally gets a short red (RGB FF3333)
# addendum run appended, matching the pattern already used elsewhere in
# the document (e.g. the "error bars" bullet near the top).

$d = $word.ActiveDocument

$wdColorGreen = 3394560   # RGB(0x00,0xCC,0x33) -> 00CC33
$wdColorRed   = 3355647   # RGB(0xFF,0x33,0x33) -> FF3333

# Unique substrings identifying each bullet that should simply turn green.
$plainTargets = @(
    "P80-81, captions of figs. 3.29-3.31",
    "P97, caption of fig.3.48",
    "P107, x-axis label of fig. 4.2",
    "P111, fig.4.6",
    "P151, label of x-axis of fig.4.50",
    "P165-168, x-axis label of figs. 5.5-5.9",
    "P188, captions of figs. 5.23 and 5.24",
    "P197-199, figs. 5.29-5.33",
    "P201-204, captions of figs. 5.33, 5.35"
)

$specialMarker = "P142-144, figs. 4.41, 4.42"

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text

    $isPlain = $false
    foreach ($marker in $plainTargets) {
        if ($t -like "*$marker*") {
            $isPlain = $true
        }
    }

    if ($isPlain) {
        $p.Range.Font.Color = $wdColorGreen
    }
    elseif ($t -like "*$specialMarker*") {
        $pRange = $p.Range

        # " into time." -> " into time. - ", scoped to this paragraph only
        # so the identical phrase elsewhere in the document is untouched.
        $findResult = $pRange.Find.Execute(" into time.", $true, $false, $false, `
            $false, $false, $true, 1, $false, " into time. - ", 2)

        # Colour the whole (now-edited) paragraph green first …
        $p.Range.Font.Color = $wdColorGreen

        # … then append the red addendum as a brand-new run after the
        # paragraph mark-minus-one (i.e. right at the end of the visible
        # text, before the pilcrow).
        $pEnd = $p.Range.End
        $insertionPoint = $d.Range($pEnd - 1, $pEnd - 1)
        $insertionPoint.InsertAfter("but would be nicer to convert axes")

        $newRunStart = $pEnd - 1
        $newRunEnd = $newRunStart + "but would be nicer to convert axes".Length
        $newRun = $d.Range($newRunStart, $newRunEnd)
        $newRun.Font.Color = $wdColorRed
    }
}

Write-Host "done"
